$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for all existing data rows (2-431)
#    from 2023-10-04 (45203) to 2023-10-05 (45204).
$ws.Range("C2:C431").Value = 45204

# 2. Row 431 regains an explicit row height (ht="15" customHeight="1").
$ws.Rows.Item(431).RowHeight = 15

# 3. Append the new record as row 432.
$ws.Range("A432").Value = "A 47440-2023"
$ws.Range("B432").Value = 45202
$ws.Range("C432").Value = 45204
$ws.Range("D432").Value = "DALARNAS LÄN"
$ws.Range("E432").Value = "MORA"
$ws.Range("F432").Value = "Bergvik skog väst AB"
$ws.Range("G432").Value = 3
$ws.Range("H432").Value = 0
$ws.Range("I432").Value = 0
$ws.Range("J432").Value = 0
$ws.Range("K432").Value = 0
$ws.Range("L432").Value = 0
$ws.Range("M432").Value = 0
$ws.Range("N432").Value = 0
$ws.Range("O432").Value = 0
$ws.Range("P432").Value = 0
$ws.Range("Q432").Value = 0

# Match formatting used by the rest of the table: B/C are dates (style with
# numFmt yyyy-mm-dd), R has wrap-text applied even though it stays empty.
$ws.Range("B432:C432").NumberFormat = "YYYY-MM-DD"
$ws.Range("R432").WrapText = $true
